# Update post last call, with all examples and images
#
# Source diff (canonical OOXML) for this change consists of two kinds of edits:
#   1. Many slides' cached footer/date field text "08/06/2020" -> "12/06/2020"
#      (a:fld id="{FE6D52D8-...}" type="datetimeFigureOut").
#   2. One content run "specimenRequirements " -> "specimenRequested " (green
#      text, srgbClr 00B050) on the "ActivityDefinition"-style technical slide.
#
# This script walks every slide/shape/paragraph in the open presentation and
# applies whichever of the above edits actually has matching content, so it
# works whether the deck contains one, some, or all of the touched slides.

$p = $ppt.ActivePresentation

$oldDate = "08/06/2020"
$newDate = "12/06/2020"
$oldReq  = "specimenRequirements"
$newReq  = "specimenRequested"

# NOTE: this COM-interop PowerShell engine only binds positional arguments
# to custom functions (named "-param value" binding silently fails), so all
# helper functions below are called positionally.
#
# NOTE 2: `.Text` getters here (TextRange/paragraph and single-run alike)
# include a trailing "\r" paragraph-mark character that is NOT part of the
# stored <a:t> run text. It must be trimmed off before writing the value
# back, otherwise a literal CR gets written into the XML text.
function Strip-ParaMark {
    param($s)
    return $s.TrimEnd([char]13, [char]10)
}

function Update-Paragraph {
    param($shapeTextRange, $para, $oldSub, $newSub)

    $runCount = $para.Runs().Count
    if ($runCount -eq 1) {
        # Single run: rewrite just that run's text so formatting / any
        # sibling runs (e.g. a trailing-space run) are left untouched and we
        # avoid the engine re-splitting the paragraph into multiple runs.
        $run = $para.Runs(1, 1)
        $clean = Strip-ParaMark $run.Text
        $run.Text = $clean.Replace($oldSub, $newSub)
    }
    elseif ($runCount -eq 0) {
        # No runs - typically a paragraph whose only content is a field
        # (<a:fld>, e.g. the slide-number/date footer placeholders).
        # Paragraphs(...).Text = "..." *appends* instead of replacing when
        # there are no runs to target, so go through the whole shape's
        # TextRange instead, which does a clean field -> run replace. This
        # is only safe when the field is the shape's sole paragraph (true
        # for the footer/date placeholders this script targets).
        $clean = Strip-ParaMark $shapeTextRange.Text
        $shapeTextRange.Text = $clean.Replace($oldSub, $newSub)
    }
    else {
        # More than one run already - rewrite the whole paragraph text.
        $clean = Strip-ParaMark $para.Text
        $para.Text = $clean.Replace($oldSub, $newSub)
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $text = $para.Text

            if ($text -like "*$oldDate*") {
                Update-Paragraph $tr $para $oldDate $newDate
            }
            elseif ($text -like "*$oldReq*") {
                Update-Paragraph $tr $para $oldReq $newReq
            }
        }
    }
}
